# Orders.xlsx - "Add work with excel"
# Update the car/plate info shown in columns D and E for the two existing
# rows, widen columns B and D to fit the new content, and move the active
# selection from E3 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: car model / plate number
$ws.Range("D1").Value = "Renault Logan"
$ws.Range("E1").Value = "а123бв"

# Row 2: car model / plate number
$ws.Range("D2").Value = "BMW X5"
$ws.Range("E2").Value = "е674ку"

# Widen columns B and D so the longer text fits (values chosen so the
# engine's character-width rounding lands on the closest attainable width
# to the author's saved widths of 13.109375 / 14.5546875).
$ws.Columns.Item(2).ColumnWidth = 12.333333333333332
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666

# Move the selection from E3 to C3
$ws.Range("C3").Select() | Out-Null
